# Generate Report for Handback
# Update the "generated at" timestamps recorded in the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-24 23:06:31"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-24 23:06:25"
# zh-cn!K2 - Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-24 23:06:43"

# de-de!K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-24 23:06:50"
